$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.9916810989379883
$ws.Range("B1").Value = 2.185224294662476
$ws.Range("C1").Value = 4.603261470794678
$ws.Range("D1").Value = 2.863547801971436
$ws.Range("E1").Value = 1.380898118019104
